$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A20").Value = "OSD slider can report timeshift buffer extent to scale, like WMC does"
$ws.Range("C20").Value = "Need to display a calculated progress based on several infotags. Can't do this in the skinning engine."

# New B20 cell - match the wrap/top-aligned style used by the rest of column B
$ws.Range("B20").Value = "Not skinnable"
$ws.Range("B20").WrapText = $true
$ws.Range("B20").VerticalAlignment = -4160

# Restore the view to where the author left it: scrolled so row 13 is at the
# top, with B21 as the active/selected cell.
$ws.Range("B21").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13

